$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format the data range as Text so purely-numeric / fraction-like
# strings (e.g. "2/2", "3", "13/13") are stored as text, matching the
# original inlineStr cell type, instead of being auto-converted to
# numbers or dates by Excel.
$ws.Range("B2:G16").NumberFormat = "@"

# Row 2
$ws.Range("B2").Value = '1210004₾'
$ws.Range("C2").Value = '29 კვ.მ.'
$ws.Range("D2").Value = '2'
$ws.Range("E2").Value = '1'
$ws.Range("F2").Value = '13/13'
$ws.Range("G2").Value = 'თბილისი, ვაკე-საბურთალო, ვაკე, ყიფშიძის ქ.'

# Row 3
$ws.Range("B3").Value = '292000₾'
$ws.Range("C3").Value = '97 კვ.მ.'
$ws.Range("D3").Value = '7/8'
$ws.Range("E3").Value = '1'
$ws.Range("F3").Value = '82'
$ws.Range("G3").Value = 'თბილისი, დიდუბე-ჩუღურეთი, დიდუბე, ა. წერეთლის გამზ.'

# Row 4
$ws.Range("B4").Value = '390000₾'
$ws.Range("C4").Value = '130 კვ.მ.'
$ws.Range("D4").Value = '4'
$ws.Range("E4").Value = '3'
$ws.Range("F4").Value = '2/2'
$ws.Range("G4").Value = 'თბილისი, დიდუბე-ჩუღურეთი, ჩუღურეთი, ი.ჯავახიშვილის ქ.'

# Row 5
$ws.Range("B5").Value = '5800924₾'
$ws.Range("C5").Value = '130 კვ.მ.'
$ws.Range("D5").Value = '2'
$ws.Range("E5").Value = '2'
$ws.Range("F5").Value = '3'
$ws.Range("G5").Value = '3/17'

# Row 6
$ws.Range("B6").Value = '4165872₾'
$ws.Range("C6").Value = '193 კვ.მ.'
$ws.Range("D6").Value = '2'
$ws.Range("E6").Value = '2'
$ws.Range("F6").Value = '6'
$ws.Range("G6").Value = '10/11'

# Row 7
$ws.Range("B7").Value = '3912252₾'
$ws.Range("C7").Value = '150 კვ.მ.'
$ws.Range("D7").Value = '5'
$ws.Range("E7").Value = '4'
$ws.Range("F7").Value = '2/3'
$ws.Range("G7").Value = 'თბილისი, ვაკე-საბურთალო, ვაკე, მცხეთის ქ.'

# Row 8
$ws.Range("B8").Value = '148500₾'
$ws.Range("C8").Value = '37 კვ.მ.'
$ws.Range("D8").Value = '2'
$ws.Range("E8").Value = '1'
$ws.Range("F8").Value = '3/12'
$ws.Range("G8").Value = 'თბილისი, ძველი თბილისი, კრწანისი, გორგასლის ქ.'

# Row 9
$ws.Range("B9").Value = '513000₾'
$ws.Range("C9").Value = '200 კვ.მ.'
$ws.Range("D9").Value = '2'
$ws.Range("E9").Value = '2'
$ws.Range("F9").Value = '4'
$ws.Range("G9").Value = '2/9'

# Row 10
$ws.Range("B10").Value = '442000₾'
$ws.Range("C10").Value = '95 კვ.მ.'
$ws.Range("D10").Value = '3'
$ws.Range("E10").Value = '2'
$ws.Range("F10").Value = '3/4'
$ws.Range("G10").Value = 'თბილისი, ვაკე-საბურთალო, საბურთალო, ბუდაპეშტის ქ.'

# Row 11
$ws.Range("B11").Value = '765000₾'
$ws.Range("C11").Value = '177 კვ.მ.'
$ws.Range("D11").Value = '5'
$ws.Range("E11").Value = '4'
$ws.Range("F11").Value = '4/4'
$ws.Range("G11").Value = 'თბილისი, ვაკე-საბურთალო, საბურთალო, მიცკევიჩის ქ.'

# Row 12
$ws.Range("B12").Value = '305000₾'
$ws.Range("C12").Value = '69 კვ.მ.'
$ws.Range("D12").Value = '2'
$ws.Range("E12").Value = '1'
$ws.Range("F12").Value = '9/14'
$ws.Range("G12").Value = 'თბილისი, ვაკე-საბურთალო, საბურთალო, ს.ცინცაძის ქ.'

# Row 13
$ws.Range("B13").Value = '156000₾'
$ws.Range("C13").Value = '61 კვ.მ.'
$ws.Range("D13").Value = '1'
$ws.Range("E13").Value = '1'
$ws.Range("F13").Value = '2'
$ws.Range("G13").Value = '3/8'

# Row 14
$ws.Range("B14").Value = '5275004₾'
$ws.Range("C14").Value = '123 კვ.მ.'
$ws.Range("D14").Value = '5'
$ws.Range("E14").Value = '3'
$ws.Range("F14").Value = '5/11'
$ws.Range("G14").Value = 'თბილისი, ვაკე-საბურთალო, საბურთალო'

# Row 15
$ws.Range("B15").Value = '2435003₾'
$ws.Range("C15").Value = '76 კვ.მ.'
$ws.Range("D15").Value = '1/5'
$ws.Range("E15").Value = '1'
$ws.Range("F15").Value = '7'
$ws.Range("G15").Value = 'თბილისი, ვაკე-საბურთალო, საბურთალო, ვაჟა-ფშაველას გამზ.'

# Row 16
$ws.Range("B16").Value = '205000₾'
$ws.Range("C16").Value = '84 კვ.მ.'
$ws.Range("D16").Value = '3'
$ws.Range("E16").Value = '2'
$ws.Range("F16").Value = '5/6'
$ws.Range("G16").Value = 'თბილისი, დიდუბე-ჩუღურეთი, დიღმის მასივი'

# Restore the default cell style (drop the temporary Text number format)
# now that the values have been written, so styling matches the original.
$ws.Range("B2:G16").Style = "Normal"
